# Password checked against API
# Applies the "insulin.xlsx" dataset refresh:
#  - renames the "Indexes.xlsx" sheet to "05-06-2021"
#  - relabels "M" -> "Male" on the General sheet
#  - updates the two patient rows' glucose/insulin readings
#  - swaps the index calculations (formulas + labels) on the Indexes sheet
#  - drops the now-unused "stumvoll1" row

$wb = $excel.ActiveWorkbook

$general = $wb.Worksheets.Item("General")
$indexes = $wb.Worksheets.Item("Indexes.xlsx")

# --- General sheet -----------------------------------------------------

$general.Range("C2").Value = "Male"

# Row 2 (first patient)
$general.Range("G2").Value = 120.0
$general.Range("H2").Value = 120.0
$general.Range("I2").Value = 110.0
$general.Range("J2").Value = 100.0
$general.Range("L2").Value = 56.0
$general.Range("M2").Value = 56.0
$general.Range("N2").Value = 55.0
$general.Range("O2").Value = 55.5
$general.Range("P2").Value = 56.0
$general.Range("R2").Value = "-"
$general.Range("T2").Value = 55.0

# Row 3 (second patient)
$general.Range("G3").Value = 6.666666666666666
$general.Range("H3").Value = 6.666666666666666
$general.Range("I3").Value = 6.11
$general.Range("J3").Value = 5.555555555555555
$general.Range("L3").Value = 336.0
$general.Range("M3").Value = 336.0
$general.Range("N3").Value = 330.0
$general.Range("O3").Value = 333.0
$general.Range("P3").Value = 336.0
$general.Range("R3").Value = "-"
$general.Range("T3").Value = 19.25

# --- Indexes sheet -------------------------------------------------------

# Row 2 : stumvoll2 -> belfiore
$indexes.Range("C2").Value = "belfiore"
$indexes.Range("D2").Formula = "=2 / (((0.5 * General!F3 + General!H3 + General!J3) / 19.08) * ((0.5 * General!L2 + General!N2 + General!P2) / 104.0) + 1 )"
$indexes.Range("E2").Value = "Healthy"
$indexes.Range("F2").Value = "≅1"

# Row 3 : cederholm -> revised
$indexes.Range("C3").Value = "revised"
$indexes.Range("D3").Formula = "=1.0 / (LN(General!F3) + LN(General!L2) + LN(General!T3))"
$indexes.Range("E3").Value = "-"
$indexes.Range("F3").Value = "0.448±0.013"

# Row 4 : matsuda -> avingon
$indexes.Range("C4").Value = "avingon"
$indexes.Range("D4").Formula = "=((0.137 * 100000000 /(General!F3 * General!L2 * 150/General!Q3)) + 100000000 /(General!J3 * General!P2 * 150/General!Q3)) / 2"
$indexes.Range("E4").Value = "-"
$indexes.Range("F4").Value = "-"

# Row 5 (stumvoll1) is retired entirely
$indexes.Rows(5).Delete()

# Column E was sized for the long "Type two diabetes"/"Insulin Resistance"
# diagnosis strings; re-fit it now that the text is short again ("Healthy"/"-").
$indexes.Columns.Item(5).ColumnWidth = 9.736979166666666

# Rename the sheet last so the lookups above (by old name) keep working.
$indexes.Name = "05-06-2021"
